{"js": "// Merge the three runs that make up each \"<id>p004r_N</id>\" paragraph\n// (\"<id>\", \"p004r_N\", \"</id>\") into a single run containing the full\n// text. Word.InsertLocation.replace keeps the formatting of the range\n// being replaced (i.e. the formatting of the leading \"<id>\" run:\n// Courier New, color 7f6000, size 9pt), so the resulting single run\n// ends up with exactly that formatting without needing to re-apply it\n// (and without perturbing any following empty run's formatting).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load('items/text');\nawait context.sync();\n\n// Identify the target paragraphs: whole paragraph text looks like\n// \"<id>...</id>\" (the tc/tcn \"<id>\" field, originally emitted as 3\n// separate runs: \"<id>\", the value, and \"</id>\").\nconst idPattern = /^<id>[\\s\\S]*<\\/id>$/;\nconst targets = [];\nfor (let i = 0; i < paragraphs.items.length; i++) {\n    const para = paragraphs.items[i];\n    if (idPattern.test(para.text)) {\n        targets.push(para);\n    }\n}\n\nfor (const para of targets) {\n    const fullText = para.text;\n    const paraRange = para.getRange();\n    // Replacing the whole paragraph range with its own text collapses\n    // the 3 runs into a single run, inheriting the first run's rPr.\n    paraRange.insertText(fullText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Merge the three runs that make up each \"<id>p004r_N</id>\" paragraph\n# (\"<id>\", \"p004r_N\", \"</id>\") into a single run containing the full\n# text. A Find/Replace over the paragraph's own (unchanged) text forces\n# Word to rebuild the paragraph as one run, which inherits the\n# formatting of the first (\"<id>\") run (Courier New, color 7f6000,\n# size 9pt) without disturbing any other run in the document.\n\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $para = $d.Paragraphs($i)\n    $rng = $para.Range\n    [void]$rng.MoveEnd(1, -1)  # wdCharacter: exclude the trailing paragraph mark\n    $text = $rng.Text\n\n    if ($text -match '^<id>[\\s\\S]*</id>$') {\n        $find = $rng.Find\n        $find.ClearFormatting()\n        $find.Replacement.ClearFormatting()\n        $find.Text = $text\n        $find.Replacement.Text = $text\n        [void]$find.Execute($text, $false, $false, $false, $false, $false, $true, 1, $false, $text, 2)\n    }\n}\n"}
